$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 46065

# Row 3
$ws.Range("A3").Value = "A 46826-2022"
$ws.Range("B3").Value = 44848
$ws.Range("C3").Value = 46065
$ws.Range("G3").Value = 4.5

# Row 4
$ws.Range("C4").Value = 46065

# Row 5
$ws.Range("A5").Value = "A 1053-2022"
$ws.Range("B5").Value = 44571
$ws.Range("C5").Value = 46065
$ws.Range("G5").Value = 1.7

# Row 6
$ws.Range("A6").Value = "A 27724-2022"
$ws.Range("B6").Value = 44743.48386574074
$ws.Range("C6").Value = 46065
$ws.Range("G6").Value = 1.3

# Row 7
$ws.Range("A7").Value = "A 64431-2023"
$ws.Range("B7").Value = 45280
$ws.Range("C7").Value = 46065
$ws.Range("G7").Value = 0.5

# Row 8
$ws.Range("A8").Value = "A 64445-2023"
$ws.Range("B8").Value = 45280
$ws.Range("C8").Value = 46065
$ws.Range("G8").Value = 3.7

# Row 9
$ws.Range("A9").Value = "A 8748-2022"
$ws.Range("B9").Value = 44613
$ws.Range("C9").Value = 46065
$ws.Range("G9").Value = 1

# Row 10
$ws.Range("A10").Value = "A 31120-2023"
$ws.Range("B10").Value = 45113
$ws.Range("C10").Value = 46065
$ws.Range("G10").Value = 0.2

# Row 11
$ws.Range("A11").Value = "A 50934-2024"
$ws.Range("B11").Value = 45602
$ws.Range("C11").Value = 46065
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("A12").Value = "A 46779-2025"
$ws.Range("B12").Value = 45926
$ws.Range("C12").Value = 46065
$ws.Range("G12").Value = 1.5

# Row 13
$ws.Range("A13").Value = "A 56917-2025"
$ws.Range("B13").Value = 45978.58453703704
$ws.Range("C13").Value = 46065
$ws.Range("G13").Value = 0.7

# Row 14
$ws.Range("A14").Value = "A 56948-2025"
$ws.Range("B14").Value = 45978.64356481482
$ws.Range("C14").Value = 46065
$ws.Range("G14").Value = 4.7
